$wb = $excel.ActiveWorkbook

# Rename the "NOTES " sheet (trailing space) to "NOTES"
$notesSheet = $wb.Worksheets.Item("NOTES ")
$notesSheet.Name = "NOTES"

# Fix up the Print_Area defined name so it refers to the new sheet name
# (it keeps the stale quoted name after a plain rename)
$wb.Names.Item("NOTES!Print_Area").RefersTo = "=NOTES!`$A`$1:`$Q`$23"

# Make NOTES the active sheet/tab (was COVER before)
$notesSheet.Activate()
$notesSheet.Select()
